# Create Credit Memo, Create Accept Payment, Create Show Address
#
# Adds the Invoice/Payment/Address columns (AB:AG) to the customer
# export sheet, refreshes the "last touched" timestamp in C2, and
# leaves the cursor on the new Amount cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- refresh the "last saved" marker for this record ------------------
$ws.Range("C2").Value = "Quick 2019/11/26 14:28:24"

# --- new header row (row 1) -------------------------------------------
$ws.Range("AB1").Value = "Invoice Number"
$ws.Range("AC1").Value = "Payment mode"
$ws.Range("AD1").Value = "Payment Mode Details"
$ws.Range("AE1").Value = "Amount"
$ws.Range("AF1").Value = "Address Name"
$ws.Range("AG1").Value = "Country"

# --- Create Credit Memo: invoice + amount details (row 2) -------------
$ws.Range("AB2").Value = "INV-523-261119-44"
$ws.Range("AB2").Style = "Normal"

# "300" must stay a text value (matches the source data export), so
# force text entry with a leading apostrophe, then drop back to the
# sheet's default (unstyled) look.
$ws.Range("AC2").Value = "'300"
$ws.Range("AC2").Style = "Normal"

# --- Create Accept Payment ---------------------------------------------
$ws.Range("AD2").Value = "Online Bank"
$ws.Range("AE2").Value = 199

# --- Create Show Address ------------------------------------------------
$ws.Range("AF2").Value = "Local"
$ws.Range("AG2").Value = "India"

# Leave the selection on the newly entered Amount cell.
$ws.Range("AE2").Select()
